$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 368, shifting existing rows 368-415 down to 369-416.
$ws.Rows(368).Insert()

# Populate the newly inserted row 368 with the new weekly price record.
$ws.Cells.Item(368, 1).Value  = 10
$ws.Cells.Item(368, 2).Value  = "Vega Modelo de Temuco"
$ws.Cells.Item(368, 3).Value  = "La Araucanía"
$ws.Cells.Item(368, 4).Value  = 45131
$ws.Cells.Item(368, 5).Value  = 9
$ws.Cells.Item(368, 6).Value  = 100112052
$ws.Cells.Item(368, 7).Value  = "Albahaca"
$ws.Cells.Item(368, 8).Value  = "Sin especificar"
$ws.Cells.Item(368, 9).Value  = "Primera"
$ws.Cells.Item(368, 10).Value = 55
$ws.Cells.Item(368, 11).Value = 6000
$ws.Cells.Item(368, 12).Value = 6000
$ws.Cells.Item(368, 13).Value = 6000
$ws.Cells.Item(368, 14).Value = "$/paquete"
$ws.Cells.Item(368, 15).Value = "Región de Arica y Parinacota"
$ws.Cells.Item(368, 16).Value = 6000
$ws.Cells.Item(368, 17).Value = 1
$ws.Cells.Item(368, 18).Value = "Hortaliza"
